# Duplicate each "congruent" trial row (column D = 1) in the condition
# table, inserting the copy directly below the original row. This adds a
# repeated congruent trial per color block (used for the speed-vs-accuracy
# csv check mentioned in the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Process from the bottom of the table upward so inserting rows doesn't
# invalidate the row numbers still to be visited.
for ($r = 10; $r -ge 2; $r--) {
    $congruent = $ws.Cells.Item($r, 4).Value()
    if ($congruent -eq 1) {
        $ws.Rows.Item($r + 1).Insert()
        $src = "A" + $r + ":D" + $r
        $dst = "A" + ($r + 1) + ":D" + ($r + 1)
        $ws.Range($src).Copy($ws.Range($dst))
    }
}

$ws.Range("H9").Select()
